# "Generate Report for Handback"
#
# The localization-status report gets refreshed once the zh-cn / de-de
# handback packages come back "in sync with en-US": the Status column
# moves from "Ready for handoff" to "Handed back: in sync with en-US",
# the new Target/Handback-File columns (E/F) get populated with links to
# the handed-back files, and the Handback DateTime (G) is stamped with
# the real timestamp instead of the 0001-01-01 placeholder.

$wb = $excel.ActiveWorkbook

function Update-HandbackSheet {
    param(
        [string]$SheetName,
        [string]$MdUrlRow2,
        [string]$XlfUrlRow2,
        [string]$MdUrlRow3,
        [string]$XlfUrlRow3,
        [string]$HandbackDateTime
    )

    $ws = $wb.Worksheets.Item($SheetName)
    Write-Host ("Updating sheet: " + $SheetName)

    $newStatus = "Handed back: in sync with en-US"

    # Row 2 (6475b226-...)
    $ws.Range("B2").Value = $newStatus

    $srcMd2 = $ws.Range("A2").Text
    $ws.Range("E2").Value = $srcMd2
    $ws.Hyperlinks.Add($ws.Range("E2"), $MdUrlRow2, "", "", $srcMd2) | Out-Null

    $srcXlf2 = $ws.Range("C2").Text
    $ws.Range("F2").Value = $srcXlf2
    $ws.Hyperlinks.Add($ws.Range("F2"), $XlfUrlRow2, "", "", $srcXlf2) | Out-Null

    $ws.Range("G2").Value = $HandbackDateTime

    # Row 3 (6d241c75-...)
    $ws.Range("B3").Value = $newStatus

    $srcMd3 = $ws.Range("A3").Text
    $ws.Range("E3").Value = $srcMd3
    $ws.Hyperlinks.Add($ws.Range("E3"), $MdUrlRow3, "", "", $srcMd3) | Out-Null

    $srcXlf3 = $ws.Range("C3").Text
    $ws.Range("F3").Value = $srcXlf3
    $ws.Hyperlinks.Add($ws.Range("F3"), $XlfUrlRow3, "", "", $srcXlf3) | Out-Null

    $ws.Range("G3").Value = $HandbackDateTime

    # Match the existing hyperlink look (underline, hyperlink blue) used by
    # the other linked cells in the sheet (A/C columns).
    $linkRange = $ws.Range("E2:F3")
    $linkRange.Font.Underline = $true
    $linkRange.Font.Color = 15570276
}

# zh-cn sheet (positional args -- named "-Param value" binding is not
# reliable in this PowerShell host, so call positionally instead)
Update-HandbackSheet `
    "zh-cn" `
    "https://github.com/OpenLocalizationTest/oltest/blob/5efcad0cd2a3cf84d5d0b84f76541ac058337007/e2e/6475b226-343f-4863-9642-7fda910e381b.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6aa48cc7a83325c93d078a8b6cfa16590fdff9fe/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/6475b226-343f-4863-9642-7fda910e381b.d14515aa11b6b7555a3b2a498308c28422409836.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTest/oltest/blob/5efcad0cd2a3cf84d5d0b84f76541ac058337007/e2e/6d241c75-725b-445d-a6f7-b30befef52dd.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6aa48cc7a83325c93d078a8b6cfa16590fdff9fe/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/6d241c75-725b-445d-a6f7-b30befef52dd.263c202ec00b6ed13443eeebc85d6fac4ee2d6e8.zh-cn.xlf" `
    "2016-03-09 00:10:51"

# de-de sheet
Update-HandbackSheet `
    "de-de" `
    "https://github.com/OpenLocalizationTest/oltest/blob/5efcad0cd2a3cf84d5d0b84f76541ac058337007/e2e/6475b226-343f-4863-9642-7fda910e381b.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5d590a855c13d424c0669e3cb57403c737dcdc94/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/6475b226-343f-4863-9642-7fda910e381b.d14515aa11b6b7555a3b2a498308c28422409836.de-de.xlf" `
    "https://github.com/OpenLocalizationTest/oltest/blob/5efcad0cd2a3cf84d5d0b84f76541ac058337007/e2e/6d241c75-725b-445d-a6f7-b30befef52dd.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5d590a855c13d424c0669e3cb57403c737dcdc94/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/6d241c75-725b-445d-a6f7-b30befef52dd.263c202ec00b6ed13443eeebc85d6fac4ee2d6e8.de-de.xlf" `
    "2016-03-09 00:11:19"

Write-Host "Handback report generated."
